# Basic team matchups added
$wb = $excel.ActiveWorkbook

$wsMatches = $wb.Worksheets.Item("swiss_matches")
$wsTeams   = $wb.Worksheets.Item("teams")

# Fix the "Cosmo" -> "Cosmoo" typo everywhere it appears (shared string
# text update), then rewrite the matchups table to two rows:
#   Cosmoo vs Frizmi            15 - 10
#   Frizmi vs Nuclear Discs     10 - 10

$wsTeams.Range("A1").Value = "Cosmoo"

$wsMatches.Range("A1").Value = "Cosmoo"
$wsMatches.Range("B1").Value = "Frizmi"
$wsMatches.Range("C1").Value = 15
$wsMatches.Range("D1").Value = 10

$wsMatches.Range("A2").Value = "Frizmi"
$wsMatches.Range("B2").Value = "Nuclear Discs"
$wsMatches.Range("C2").Value = 10
$wsMatches.Range("D2").Value = 10

# Remove the old third matchup row entirely
$wsMatches.Range("A3:D3").Delete()

# Match the author's final selection / active cell
$wsMatches.Range("A2").Select()
